$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark wherever it currently sits
#    (it is currently wrapped around the title text, right after "MP73010").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Find the paragraph that holds "Ben changing things up!" and insert a
#    brand-new paragraph straight after it containing "I love Programming".
$newPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Ben changing things up!*") {
        $para.Range.InsertParagraphAfter()
        $newPara = $para.Next()
        # Type a trailing sentinel character along with the real text so the
        # bookmark-anchor position we need (immediately after "Programming",
        # before the paragraph mark) is not the literal last slot of the
        # paragraph when we create the bookmark.
        $newPara.Range.Text = "I love ProgrammingX"
        break
    }
}

# 3. Re-add the "_GoBack" bookmark right after "I love Programming" (and
#    before the sentinel character), then drop the sentinel so the bookmark
#    ends up collapsed immediately before the paragraph mark, matching the
#    target document.
$pos = $newPara.Range.End - 2
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$sentinel = $d.Range($newPara.Range.End - 2, $newPara.Range.End - 1)
$sentinel.Delete()
